$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply Text number format to D2:E51 so numeric-looking strings (e.g. "1.000")
# are preserved exactly as text instead of being auto-converted to numbers by Excel.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '27.584.68'
$ws.Range("E2").Value = '  -2.58%  '
$ws.Range("D3").Value = '1.753.58'
$ws.Range("E3").Value = '  -3.41%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '324.45'
$ws.Range("E5").Value = '  -0.30%  '
$ws.Range("D6").Value = '1.000'
$ws.Range("D7").Value = '0.4487'
$ws.Range("E7").Value = '  +3.19%  '
$ws.Range("D8").Value = '0.3619'
$ws.Range("E8").Value = '  -1.33%  '
$ws.Range("D9").Value = '0.07510'
$ws.Range("E9").Value = '  -1.97%  '
$ws.Range("D10").Value = '42.12'
$ws.Range("E10").Value = '  -6.27%  '
$ws.Range("D11").Value = '1.105'
$ws.Range("E11").Value = '  -3.43%  '
$ws.Range("E12").Value = '  +0.03%  '
$ws.Range("D13").Value = '20.72'
$ws.Range("E13").Value = '  -5.81%  '
$ws.Range("D14").Value = '6.045'
$ws.Range("E14").Value = '  -4.37%  '
$ws.Range("E15").Value = '  -3.92%  '
$ws.Range("D16").Value = '1.750.54'
$ws.Range("E16").Value = '  -3.78%  '
$ws.Range("D17").Value = '92.96'
$ws.Range("E17").Value = '  -2.57%  '
$ws.Range("E18").Value = '  -1.30%  '
$ws.Range("D19").Value = '0.06395'
$ws.Range("E19").Value = '  -0.69%  '
$ws.Range("D20").Value = '0.9999'
$ws.Range("E21").Value = '  -2.95%  '
$ws.Range("D22").Value = '5.854'
$ws.Range("E22").Value = '  -6.09%  '
$ws.Range("D23").Value = '27.617.43'
$ws.Range("E23").Value = '  -2.49%  '
$ws.Range("E24").Value = '  -2.95%  '
$ws.Range("D25").Value = '2.108'
$ws.Range("E25").Value = '  -1.44%  '
$ws.Range("D26").Value = '161.84'
$ws.Range("E26").Value = '  +1.07%  '
$ws.Range("D27").Value = '20.38'
$ws.Range("E27").Value = '  -1.63%  '
$ws.Range("D28").Value = '1.953.36'
$ws.Range("E28").Value = '  -3.78%  '
$ws.Range("D29").Value = '2.123'
$ws.Range("E29").Value = '  -6.63%  '
$ws.Range("D30").Value = '125.38'
$ws.Range("E30").Value = '  -4.08%  '
$ws.Range("D31").Value = '1.082'
$ws.Range("E31").Value = '  -9.97%  '
$ws.Range("B32").Value = 'HuobiToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D32").Value = '3.658'
$ws.Range("E32").Value = '  +2.83%  '
$ws.Range("B33").Value = 'Stellar'
$ws.Range("C33").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D33").Value = '0.09031'
$ws.Range("E33").Value = '  -1.05%  '
$ws.Range("D34").Value = '5.556'
$ws.Range("E34").Value = '  -7.69%  '
$ws.Range("E35").Value = '  -7.84%  '
$ws.Range("E36").Value = '  -3.33%  '
$ws.Range("B37").Value = 'Algorand'
$ws.Range("C37").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D37").Value = '0.2096'
$ws.Range("E37").Value = '  -3.53%  '
$ws.Range("B38").Value = 'TheSandbox'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D38").Value = '0.6388'
$ws.Range("E38").Value = '  -2.98%  '
$ws.Range("D39").Value = '4.991'
$ws.Range("E39").Value = '  -4.41%  '
$ws.Range("D40").Value = '0.05983'
$ws.Range("E40").Value = '  -3.45%  '
$ws.Range("D41").Value = '1.196'
$ws.Range("E41").Value = '  -0.23%  '
$ws.Range("D42").Value = '0.9994'
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("D43").Value = '1.385'
$ws.Range("E43").Value = '  -3.23%  '
$ws.Range("E44").Value = '  -3.00%  '
$ws.Range("D45").Value = '13.27'
$ws.Range("E45").Value = '  -4.39%  '
$ws.Range("D46").Value = '0.5910'
$ws.Range("E46").Value = '  -3.08%  '
$ws.Range("D47").Value = '3.712'
$ws.Range("E47").Value = '  -0.58%  '
$ws.Range("E48").Value = '  -2.80%  '
$ws.Range("D49").Value = '121.77'
$ws.Range("E49").Value = '  -3.01%  '
$ws.Range("E50").Value = '  -0.35%  '
$ws.Range("D51").Value = '0.06871'
$ws.Range("E51").Value = '  -1.81%  '

# Restore the default "Normal" style on the data range so no stray number-format
# style remains attached to the cells (matches original styling).
$dataRange.Style = "Normal"

Write-Host "Applied cryptos list update"